# Update the "Level 2" evaluation-metric block (rows 37-39) on Sheet1 with
# the new clustering run results ("evaluation metric for clustering algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header tweak: "Outliers" -> "New Outliers" for the Level 2 table.
$ws.Range("I38").Value = "New Outliers"

# Updated results for the HDBSCAN row (row 39) under "Level 2".
$ws.Range("C39").Value = "309 clusters"
$ws.Range("D39").Value = "86 outliers"
$ws.Range("F39").Value = "244.744 secs"
$ws.Range("H39").Value = 80
$ws.Range("I39").Value = 39
$ws.Range("K39").Value = "54.398 secs"

# E39 and J39 hold numeric-looking text ("0.151" / "-0.055") that must stay
# text (like the original "0.141" / "-0.095"), not get auto-converted to a
# number. Writing directly into those cells through .Value auto-coerces
# numeric-looking strings into real numbers, so:
#  - J39 already carries a quoted-text style, so prefixing with an
#    apostrophe is enough to force text while keeping its own style.
$ws.Range("J39").Value = "'-0.055"

#  - E39's style is a plain (non quote-prefixed) style, so force the new
#    value through a scratch cell typed as text, then paste *values only*
#    into E39 - this keeps E39's original style index untouched while
#    swapping in the text value.
$scratch = $ws.Range("Z1")
$scratch.Value = "'0.151"
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()

# Put the active selection on C39, matching the workbook state after the edit.
$ws.Range("C39").Select()
